$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.011.91"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "1.789.03"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").Value = "  -0.11%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "227.53"
$r.NumberFormat = "General"
$ws.Range("E5").Value = "  -1.67%  "
$ws.Range("E6").Value = "  +1.25%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  -1.97%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "46.05"
$r.NumberFormat = "General"
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("E11").Value = "  -2.94%  "
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").Value = "2.047.01"
$ws.Range("E13").Value = "  -1.58%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "11.36"
$r.NumberFormat = "General"
$ws.Range("E14").Value = "  +9.78%  "
$ws.Range("D15").Value = "1.797.38"
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("E16").Value = "  -1.69%  "
$ws.Range("D17").Value = "34.046.72"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("E18").Value = "  -2.89%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "69.51"
$r.NumberFormat = "General"
$ws.Range("E19").Value = "  -1.73%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "252.95"
$r.NumberFormat = "General"
$ws.Range("E20").Value = "  -3.05%  "
$ws.Range("D21").Value = "0.0₃0743"
$ws.Range("E21").Value = "  -1.07%  "
$ws.Range("E22").Value = "  -0.36%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "10.39"
$r.NumberFormat = "General"
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("E24").Value = "  -2.80%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "2.14"
$r.NumberFormat = "General"
$ws.Range("E25").Value = "  -2.69%  "
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "157.29"
$r.NumberFormat = "General"
$ws.Range("E26").Value = "  -2.56%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "16.57"
$r.NumberFormat = "General"
$ws.Range("E27").Value = "  -1.18%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "7.01"
$r.NumberFormat = "General"
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("E29").Value = "  -2.56%  "
$ws.Range("E30").Value = "  -0.09%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "3.85"
$r.NumberFormat = "General"
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("E32").Value = "  -0.45%  "
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "1.21"
$r.NumberFormat = "General"
$ws.Range("E33").Value = "  -1.01%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "3.61"
$r.NumberFormat = "General"
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("D36").Value = "1.463.41"
$ws.Range("E36").Value = "  -8.13%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "1.07"
$r.NumberFormat = "General"
$ws.Range("E37").Value = "  +0.87%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "0.631"
$r.NumberFormat = "General"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("E39").Value = "  -1.46%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "83.36"
$r.NumberFormat = "General"
$ws.Range("E40").Value = "  -3.02%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "2.84"
$r.NumberFormat = "General"
$ws.Range("E41").Value = "  +2.51%  "
$ws.Range("E42").Value = "  -0.47%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "0.900"
$r.NumberFormat = "General"
$ws.Range("E43").Value = "  -2.38%  "
$ws.Range("E44").Value = "  -2.96%  "
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("E46").Value = "  +0.43%  "
$ws.Range("D47").Value = "1.947.27"
$ws.Range("E47").Value = "  -1.15%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "5.71"
$r.NumberFormat = "General"
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.NumberFormat = "General"
$ws.Range("E49").Value = "  -0.02%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "11.72"
$r.NumberFormat = "General"
$ws.Range("E50").Value = "  +3.13%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "51.16"
$r.NumberFormat = "General"
$ws.Range("E51").Value = "  -4.73%  "
